$d = $word.ActiveDocument

# Replace the "address" field reference with "one_line"
$d.Content.Find.Execute("users[0].address.address", $true, $false, $false, $false, $false,
                         $true, 1, $false, "users[0].address.one_line", 2)

# Replace the "phone_number" field reference with "mobile_number"
$d.Content.Find.Execute("users[0].phone_number", $true, $false, $false, $false, $false,
                         $true, 1, $false, "users[0].mobile_number", 2)
